$d = $word.ActiveDocument

$replacements = @(
    @("83×61=5063", "46×61=2806"),
    @("61×50=3050", "28×26=728"),
    @("29×52=1508", "84×20=1680"),
    @("97×13=1261", "47×25=1175"),
    @("61×58=3538", "63×74=4662"),
    @("17×50=850", "48×14=672"),
    @("59×75=4425", "98×60=5880"),
    @("26×90=2340", "18×86=1548"),
    @("42×38=1596", "39×28=1092"),
    @("53×14=742", "63×83=5229"),
    @("52×66=3432", "15×38=570"),
    @("25×48=1200", "82×18=1476"),
    @("26×34=884", "86×84=7224"),
    @("98×92=9016", "45×38=1710"),
    @("45×85=3825", "27×70=1890"),
    @("51×44=2244", "23×93=2139"),
    @("75×61=4575", "42×37=1554"),
    @("34×11=374", "26×45=1170"),
    @("62×36=2232", "50×39=1950"),
    @("85×90=7650", "60×92=5520"),
    @("60×26=1560", "34×43=1462"),
    @("51×90=4590", "72×20=1440"),
    @("12×53=636", "53×44=2332"),
    @("28×98=2744", "60×97=5820"),
    @("16×92=1472", "21×51=1071")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
